# Commit: "Add front page & Add "E"s to titles"
#
# 1) Refresh the cached "date last edited" placeholder text that lives on
#    the slide master and every slide layout (01/01/2024 -> 02/01/2024).
# 2) Append an "E" to the "Exercise 10" title text on every slide so it
#    reads "Exercise E10".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "01/01/2024") {
                $tr.Text = "02/01/2024"
            }
        }
    }
}

# --- Slide master date placeholder ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- Every custom (slide) layout's date placeholder ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- Title text on every slide: "Exercise 10" -> "Exercise E10" ---
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "Exercise 10") {
                $tr.Text = "Exercise E10"
            }
        }
    }
}
